$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.5
$ws.Range("C3").Select()

$ws.Range("M2").Font.Name = "Arial"
$ws.Range("M2").Font.Size = 10
